# Update the locations sheet: header stays the same, but the three data
# rows now hold real-world SJSU-area landmarks instead of the placeholder
# "Location N" rows, with refreshed lat/long coordinates.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (unchanged text, but re-asserted for safety)
$ws.Range("A1").Value = "Location"
$ws.Range("B1").Value = "Latitude"
$ws.Range("C1").Value = "Longitude"

# Row 3 is entered before row 2 so that the shared-string table ends up in
# the same order as the authored workbook (Martin Luther King Jr Library
# before San Jose State University).
$ws.Range("A3").Value = "Martin Luther King Jr Library"
$ws.Range("B3").Value = 37.3355
$ws.Range("C3").Value = -121.885

$ws.Range("A2").Value = "San Jose State University"
$ws.Range("B2").Value = 37.3352
$ws.Range("C2").Value = -121.8811

$ws.Range("A4").Value = "8th & San Fernando"
$ws.Range("B4").Value = 37.33357
$ws.Range("C4").Value = -121.87859

# Column A is widened to fit the longer location names.
$ws.Columns.Item(1).ColumnWidth = 22.17

# Leave the selection on C5, just below the data, as in the saved file.
$ws.Range("C5").Select()
